# Reorders the SectorGroup sheet's D/E/F/G columns.
#
# Before:  D=codeforiati:group-name   E=codeforiati:group-code
#          F=codeforiati:category-code G=codeforiati:category-name
# After:   D=codeforiati:group-code   E=codeforiati:category-name
#          F=codeforiati:group-name   G=codeforiati:category-code
#
# i.e. for every row: new_D = old_E, new_E = old_G, new_F = old_D, new_G = old_F

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $oldD = $ws.Cells.Item($r, 4).Value()
    $oldE = $ws.Cells.Item($r, 5).Value()
    $oldF = $ws.Cells.Item($r, 6).Value()
    $oldG = $ws.Cells.Item($r, 7).Value()

    $ws.Cells.Item($r, 4).Value = $oldE
    $ws.Cells.Item($r, 5).Value = $oldG
    $ws.Cells.Item($r, 6).Value = $oldD
    $ws.Cells.Item($r, 7).Value = $oldF
}
